# Auto-generated edit script applying scheduled market-price refresh updates
# to the Leve profit-tracking sheets (columns H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 263
$ws.Range("I4").Value = 263
$ws.Range("K4").Value = 263
$ws.Range("M4").Value = -149
$ws.Range("H18").Value = 4999
$ws.Range("I18").Value = 4999
$ws.Range("K18").Value = 4999
$ws.Range("M18").Value = -4715
$ws.Range("H19").Value = 3508
$ws.Range("I19").Value = 2807
$ws.Range("K19").Value = 2807
$ws.Range("M19").Value = -2632
$ws.Range("H40").Value = 4497.5
$ws.Range("J40").Value = 4995
$ws.Range("L40").Value = 4995
$ws.Range("N40").Value = -5345
$ws.Range("H51").Value = 37997.8
$ws.Range("J51").Value = 6663
$ws.Range("L51").Value = 6663
$ws.Range("N51").Value = -7631
$ws.Range("H70").Value = 3975.125
$ws.Range("I70").Value = 2566.6667
$ws.Range("K70").Value = 7700.000100000001
$ws.Range("M70").Value = -7430.000100000001
$ws.Range("H73").Value = 3975.125
$ws.Range("I73").Value = 2566.6667
$ws.Range("K73").Value = 7700.000100000001
$ws.Range("M73").Value = -6764.000100000001
$ws.Range("H98").Value = 2915.375
$ws.Range("I98").Value = 2762.25
$ws.Range("K98").Value = 2762.25
$ws.Range("M98").Value = -1264.25
$ws.Range("H122").Value = 2915.375
$ws.Range("I122").Value = 2762.25
$ws.Range("K122").Value = 8286.75
$ws.Range("M122").Value = -5836.75
$ws.Range("H135").Value = 496
$ws.Range("I135").Value = 496
$ws.Range("K135").Value = 4464
$ws.Range("M135").Value = -1929
$ws.Range("H138").Value = 2484.125
$ws.Range("J138").Value = 3001.9656
$ws.Range("L138").Value = 9005.8968
$ws.Range("N138").Value = -19285.8968

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8929.647000000001
$ws.Range("J2").Value = 2533.111
$ws.Range("L2").Value = 2533.111
$ws.Range("N2").Value = -2759.111
$ws.Range("H32").Value = 1828.6451
$ws.Range("I32").Value = 1828.6451
$ws.Range("K32").Value = 1828.6451
$ws.Range("M32").Value = -1541.6451
$ws.Range("H45").Value = 2499.75
$ws.Range("I45").Value = 2666.3333
$ws.Range("K45").Value = 2666.3333
$ws.Range("M45").Value = -2289.3333
$ws.Range("H110").Value = 4888.9414
$ws.Range("I110").Value = 4354.7144
$ws.Range("K110").Value = 4354.7144
$ws.Range("M110").Value = -2309.7144
$ws.Range("H116").Value = 8929.647000000001
$ws.Range("J116").Value = 2533.111
$ws.Range("L116").Value = 2533.111
$ws.Range("N116").Value = -7121.111
$ws.Range("H122").Value = 1901.1111
$ws.Range("J122").Value = 1270
$ws.Range("L122").Value = 3810
$ws.Range("N122").Value = -8710

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8929.647000000001
$ws.Range("J3").Value = 2533.111
$ws.Range("L3").Value = 2533.111
$ws.Range("N3").Value = -2761.111
$ws.Range("H24").Value = 4485.6665
$ws.Range("I24").Value = 1379.2
$ws.Range("K24").Value = 1379.2
$ws.Range("M24").Value = -1144.2
$ws.Range("H34").Value = 29999.5
$ws.Range("I34").Value = 29999
$ws.Range("K34").Value = 29999
$ws.Range("M34").Value = -29885
$ws.Range("H134").Value = 62524164
$ws.Range("I134").Value = 62524164
$ws.Range("K134").Value = 187572492
$ws.Range("M134").Value = -187569957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1187.0769
$ws.Range("I31").Value = 1010.5
$ws.Range("J31").Value = 1469.6
$ws.Range("K31").Value = 1010.5
$ws.Range("L31").Value = 1469.6
$ws.Range("M31").Value = -715.5
$ws.Range("N31").Value = -2059.6
$ws.Range("H34").Value = 1187.0769
$ws.Range("I34").Value = 1010.5
$ws.Range("J34").Value = 1469.6
$ws.Range("K34").Value = 1010.5
$ws.Range("L34").Value = 1469.6
$ws.Range("M34").Value = -808.5
$ws.Range("N34").Value = -1873.6
$ws.Range("H58").Value = 3695.0833
$ws.Range("I58").Value = 3713.1667
$ws.Range("J58").Value = 3677
$ws.Range("K58").Value = 3713.1667
$ws.Range("L58").Value = 3677
$ws.Range("M58").Value = -3510.1667
$ws.Range("N58").Value = -4083
$ws.Range("H136").Value = 3695.0833
$ws.Range("I136").Value = 3713.1667
$ws.Range("J136").Value = 3677
$ws.Range("K136").Value = 11139.5001
$ws.Range("L136").Value = 11031
$ws.Range("M136").Value = -8589.500100000001
$ws.Range("N136").Value = -16131

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 141.3125
$ws.Range("J2").Value = 180.5
$ws.Range("L2").Value = 1083
$ws.Range("N2").Value = -1309
$ws.Range("H38").Value = 70.703705
$ws.Range("I38").Value = 21.333334
$ws.Range("K38").Value = 64.00000199999999
$ws.Range("M38").Value = 282.999998
$ws.Range("H80").Value = 4618.1
$ws.Range("I80").Value = 4399.3335
$ws.Range("K80").Value = 13198.0005
$ws.Range("M80").Value = -12262.0005
$ws.Range("H83").Value = 4618.1
$ws.Range("I83").Value = 4399.3335
$ws.Range("K83").Value = 39594.0015
$ws.Range("M83").Value = -34914.0015
$ws.Range("H113").Value = 377.25
$ws.Range("J113").Value = 544.8333
$ws.Range("L113").Value = 1634.4999
$ws.Range("N113").Value = -5974.4999
$ws.Range("H122").Value = 47120.184
$ws.Range("I122").Value = 955.4
$ws.Range("K122").Value = 8598.6
$ws.Range("M122").Value = -6148.6
$ws.Range("H131").Value = 350895.12
$ws.Range("I131").Value = 657.9
$ws.Range("J131").Value = 434284.94
$ws.Range("K131").Value = 1973.7
$ws.Range("L131").Value = 1302854.82
$ws.Range("M131").Value = 3066.3
$ws.Range("N131").Value = -1312934.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4228601
$ws.Range("I11").Value = 4228601
$ws.Range("K11").Value = 4228601
$ws.Range("M11").Value = -4228462
$ws.Range("H20").Value = 3353999.8
$ws.Range("J20").Value = 30999.5
$ws.Range("L20").Value = 30999.5
$ws.Range("N20").Value = -31489.5
$ws.Range("H122").Value = 3856.5334
$ws.Range("I122").Value = 4603.4
$ws.Range("J122").Value = 2362.8
$ws.Range("K122").Value = 13810.2
$ws.Range("L122").Value = 7088.400000000001
$ws.Range("M122").Value = -11360.2
$ws.Range("N122").Value = -11988.4
$ws.Range("H123").Value = 46661.668
$ws.Range("J123").Value = 46661.668
$ws.Range("L123").Value = 46661.668
$ws.Range("N123").Value = -51561.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 25042.857
$ws.Range("J3").Value = 25042.857
$ws.Range("L3").Value = 25042.857
$ws.Range("N3").Value = -25266.857
$ws.Range("H15").Value = 25042.857
$ws.Range("J15").Value = 25042.857
$ws.Range("L15").Value = 25042.857
$ws.Range("N15").Value = -25382.857
$ws.Range("H16").Value = 2067.4285
$ws.Range("I16").Value = 2179.4167
$ws.Range("J16").Value = 1395.5
$ws.Range("K16").Value = 2179.4167
$ws.Range("L16").Value = 1395.5
$ws.Range("M16").Value = -2009.4167
$ws.Range("N16").Value = -1735.5
$ws.Range("H20").Value = 15679.8
$ws.Range("J20").Value = 15679.8
$ws.Range("L20").Value = 15679.8
$ws.Range("N20").Value = -16131.8
$ws.Range("H40").Value = 5393.6
$ws.Range("I40").Value = 3539.8
$ws.Range("J40").Value = 7247.4
$ws.Range("K40").Value = 3539.8
$ws.Range("L40").Value = 7247.4
$ws.Range("M40").Value = -3403.8
$ws.Range("N40").Value = -7519.4
$ws.Range("H55").Value = 1721.7778
$ws.Range("I55").Value = 1665.6666
$ws.Range("J55").Value = 1749.8334
$ws.Range("K55").Value = 1665.6666
$ws.Range("L55").Value = 1749.8334
$ws.Range("M55").Value = -1492.6666
$ws.Range("N55").Value = -2095.8334
$ws.Range("H68").Value = 6422.8
$ws.Range("J68").Value = 11501.25
$ws.Range("L68").Value = 11501.25
$ws.Range("N68").Value = -12999.25
$ws.Range("H71").Value = 6422.8
$ws.Range("J71").Value = 11501.25
$ws.Range("L71").Value = 57506.25
$ws.Range("N71").Value = -64994.25
$ws.Range("H136").Value = 4551.852
$ws.Range("I136").Value = 3263.0527
$ws.Range("J136").Value = 7612.75
$ws.Range("K136").Value = 9789.158100000001
$ws.Range("L136").Value = 22838.25
$ws.Range("M136").Value = -7239.158100000001
$ws.Range("N136").Value = -27938.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 26799.285
$ws.Range("J81").Value = 97500
$ws.Range("L81").Value = 195000
$ws.Range("N81").Value = -197122
$ws.Range("H84").Value = 26799.285
$ws.Range("J84").Value = 97500
$ws.Range("L84").Value = 975000
$ws.Range("N84").Value = -985608
$ws.Range("H122").Value = 2938.9375
$ws.Range("I122").Value = 3372.8
$ws.Range("K122").Value = 10118.4
$ws.Range("M122").Value = -7668.400000000001
